$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "manos" worksheet after the last existing sheet
#    ("jefe vampiro") so it becomes sheet11 / the new active tab.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "manos"

# ---------------------------------------------------------------------------
# 2. Populate the sprite data (binary strings in columns A/B/E/F, hex
#    conversion formulas in columns C/D/G/H). Row 17 is intentionally left
#    blank, matching the source sheet.
# ---------------------------------------------------------------------------
$data = @(
    @(1,"n10000001","n00000000","n00000000","n00000000"),
    @(2,"n01000001","n11000001","n00000000","n00000000"),
    @(3,"n00110000","n11100001","n00000000","n00000000"),
    @(4,"n00011100","n01110001","n00000000","n00000000"),
    @(5,"n00001110","n01110001","n00000000","n00000000"),
    @(6,"n00000111","n00110011","n00000111","n00110000"),
    @(7,"n00000011","n00111011","n00000111","n00111011"),
    @(8,"n00000011","n11111011","n00000111","n11111111"),
    @(9,"n11111111","n11111011","n00011111","n11111011"),
    @(10,"n01111111","n11011111","n00111111","n11011111"),
    @(11,"n00000011","n10111111","n00111011","n11111111"),
    @(12,"n00011111","n01111111","n00011111","n01111111"),
    @(13,"n00111101","n11111111","n00011111","n11111111"),
    @(14,"n00110001","n11111111","n00001101","n11111111"),
    @(15,"n00100000","n11111110","n00001110","n11111110"),
    @(16,"n00100000","n01111100","n00000000","n01111100"),
    @(18,"n00000000","n10000001","n00000000","n00000000"),
    @(19,"n10000011","n10000010","n00000000","n00000000"),
    @(20,"n10000111","n00001100","n00000000","n00000000"),
    @(21,"n10001110","n00111000","n00000000","n00000000"),
    @(22,"n10001110","n01110000","n00000000","n00000000"),
    @(23,"n11001100","n11100000","n00001100","n11100000"),
    @(24,"n11011100","n11000000","n11011100","n11100000"),
    @(25,"n11011111","n11000000","n11111111","n11100000"),
    @(26,"n11011111","n11111111","n11011111","n11111000"),
    @(27,"n11111011","n11111110","n11111011","n11111100"),
    @(28,"n11111101","n11000000","n11111111","n11011100"),
    @(29,"n11111110","n11111000","n11111110","n11111000"),
    @(30,"n11111111","n10111100","n11111111","n11111000"),
    @(31,"n11111111","n10001100","n11111111","n10110000"),
    @(32,"n01111111","n00000100","n01111111","n01110000"),
    @(33,"n00111110","n00000100","n00111110","n00000000")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r,1).Value = $row[1]
    $ws.Cells.Item($r,2).Value = $row[2]
    $ws.Cells.Item($r,3).Formula = '=CONCATENATE("#",BIN2HEX(REPLACE(A' + $r + ',1,1,""),2))'
    $ws.Cells.Item($r,4).Formula = '=CONCATENATE("#",BIN2HEX(REPLACE(B' + $r + ',1,1,""),2))'
    $ws.Cells.Item($r,5).Value = $row[3]
    $ws.Cells.Item($r,6).Value = $row[4]
    $ws.Cells.Item($r,7).Formula = '=CONCATENATE("#",BIN2HEX(REPLACE(E' + $r + ',1,1,""),2))'
    $ws.Cells.Item($r,8).Formula = '=CONCATENATE("#",BIN2HEX(REPLACE(F' + $r + ',1,1,""),2))'
}

# ---------------------------------------------------------------------------
# 3. View state: "manos" becomes the selected/active tab with G17 selected;
#    the previously active sheet ("jefe vampiro") drops its range selection
#    down to a single cell C1.
# ---------------------------------------------------------------------------
$ws10 = $wb.Worksheets.Item(10)
$ws10.Range("C1").Select()

$ws.Range("G17").Select()
